$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 401, pushing the existing rows 401-421 down to 402-422
$ws.Rows.Item(401).Insert()

# Populate the newly inserted row 401 with the new weekly record
$ws.Range("A401").Value = 3
$ws.Range("B401").Value = "Femacal de La Calera"
$ws.Range("C401").Value = "Coquimbo"
$ws.Range("D401").Value = 45267
$ws.Range("E401").Value = 5
$ws.Range("F401").Value = "Fruta"
$ws.Range("G401").Value = 100101
$ws.Range("H401").Value = "Berries"
$ws.Range("I401").Value = 100101001
$ws.Range("J401").Value = "Arándano (blue)"
$ws.Range("K401").Value = "Sin especificar"
$ws.Range("L401").Value = "Primera"
$ws.Range("M401").Value = 40
$ws.Range("N401").Value = 6000
$ws.Range("O401").Value = 6000
$ws.Range("P401").Value = 6000
$ws.Range("Q401").Value = "$/bandeja 2 kilos"
$ws.Range("R401").Value = "Provincia de Quillota"
$ws.Range("S401").Value = 3000
$ws.Range("T401").Value = 2
